$d = $word.ActiveDocument

$target = $d.Content
$target.Find.Execute("On startup- introduce Ms. Minutes type assistant", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$p1 = $target.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Long distance relationship"
